$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 0.1531268371546149
$ws.Range("G2").Value = 0.5341187720458554
$ws.Range("H2").Value = 2.073148148148148

$ws.Range("E3").Value = 1.310344827586207
$ws.Range("F3").Value = 0.6425720956486043
$ws.Range("G3").Value = 1.147205972906404
$ws.Range("H3").Value = 3.752241379310345

$ws.Range("E4").Value = 1.12280701754386
$ws.Range("F4").Value = 0.405733761487051
$ws.Range("G4").Value = 1.123444688109162
$ws.Range("H4").Value = 3.481929824561404

$ws.Range("F5").Value = 0.2005680555555555
$ws.Range("G5").Value = 0.6531014770723105
$ws.Range("H5").Value = 2.411333333333333

$ws.Range("E6").Value = 1.388888888888889
$ws.Range("F6").Value = 0.7214347075249853
$ws.Range("G6").Value = 1.679037716784245
$ws.Range("H6").Value = 4.601296296296296

$ws.Range("E7").Value = 1.122448979591837
$ws.Range("F7").Value = 0.4891171849692258
$ws.Range("G7").Value = 1.300359511661808
$ws.Range("H7").Value = 3.898979591836734

$ws.Range("F8").Value = 0.1155166761148904
$ws.Range("G8").Value = 0.7168164682539683
$ws.Range("H8").Value = 2.202380952380953

$ws.Range("E9").Value = 1.1875
$ws.Range("F9").Value = 0.5104406156994048
$ws.Range("G9").Value = 1.902392578125
$ws.Range("H9").Value = 5.13953125

$ws.Range("E10").Value = 1.063492063492063
$ws.Range("F10").Value = 0.3262371346686823
$ws.Range("G10").Value = 1.726775714915596
$ws.Range("H10").Value = 4.581904761904762
